$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) date column (C) for the data rows 2-41
# needs to move forward by one day (45243 -> 45244), reflecting an
# automatic daily refresh of the source data. The date is stored as a
# serial number with an existing date format (style), so we just update
# the underlying value.
$ws.Range("C2:C41").Value = 45244
